# Check validity of equilibrium equations using PE, 0.5PP, PP and 2PP cases
# -> the elastic modulus E (column AO, rows 3-18) is updated to the new
#    validated value; the dependent mi/psi-bar/i and sum columns (AR:AU)
#    recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO3:AO18").Value = 28798.820050585658

# Leave the selection on the first cell of the edited range, matching the
# single-cell selection left behind after the edit.
$ws.Range("AO3").Select()
